# Automated data refresh: the DCX row (ticker row 2) has dropped out of the
# scrape entirely, so HCTI/SER/XHLD shift up by one row. The institutional
# ownership figures (column F) for the surviving tickers were also
# re-scraped with slightly different values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the DCX row - this shifts HCTI/SER/XHLD up (rows 3,4,5 -> 2,3,4)
# and Excel automatically shrinks the used-range dimension to A1:AU4.
$ws.Rows(2).Delete()

# Refresh the "Institutional Ownership" (column F) values for the three
# remaining tickers with the newly scraped figures.
$ws.Range("F2").Value = 0.058819998
$ws.Range("F3").Value = 0.04325
$ws.Range("F4").Value = 0.0064
